$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the TPM-derived metrics for the Fbn1-Itgb3 LR-pair sheet
# (recomputed with new TPM values; only cells that changed are set)

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.445154666666667
$ws.Range("H2").Value = 10.335464
$ws.Range("I2").Value = 0.01110365039942287
$ws.Range("J2").Value = 0.01110365039942286
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 0.3808687387093334
$ws.Range("R2").Value = 3.427818648384
$ws.Range("S2").Value = 0.0001250548139589593
$ws.Range("T2").Value = 0.0001250548139589593
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.445154666666667
$ws.Range("H3").Value = 10.335464
$ws.Range("I3").Value = 0.01110365039942287
$ws.Range("J3").Value = 0.01110365039942286
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 31.04797731359644
$ws.Range("R3").Value = 279.431795822368
$ws.Range("S3").Value = 0.01019432321988742
$ws.Range("T3").Value = 0.01019432321988742
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.445154666666667
$ws.Range("H4").Value = 10.335464
$ws.Range("I4").Value = 0.01110365039942287
$ws.Range("J4").Value = 0.01110365039942286
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 2.388591188338666
$ws.Range("R4").Value = 21.497320695048
$ws.Range("S4").Value = 0.000784272365576486
$ws.Range("T4").Value = 0.000784272365576486
$ws.Range("I5").Value = 0.938949437922138
$ws.Range("J5").Value = 0.938949437922138
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 32.20710984847199
$ws.Range("R5").Value = 289.863988636248
$ws.Range("S5").Value = 0.010574913929416
$ws.Range("T5").Value = 0.010574913929416
$ws.Range("I6").Value = 0.938949437922138
$ws.Range("J6").Value = 0.938949437922138
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("S6").Value = 0.8620547039023684
$ws.Range("T6").Value = 0.8620547039023685
$ws.Range("I7").Value = 0.938949437922138
$ws.Range("J7").Value = 0.938949437922138
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 201.984597230559
$ws.Range("R7").Value = 1817.861375075031
$ws.Range("S7").Value = 0.06631982009035357
$ws.Range("T7").Value = 0.06631982009035359
$ws.Range("G8").Value = 15.497141
$ws.Range("H8").Value = 46.491423
$ws.Range("I8").Value = 0.04994691167843914
$ws.Range("J8").Value = 0.04994691167843914
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 1.713239931832
$ws.Range("R8").Value = 15.419159386488
$ws.Range("S8").Value = 0.0005625268738735173
$ws.Range("T8").Value = 0.0005625268738735174
$ws.Range("G9").Value = 15.497141
$ws.Range("H9").Value = 46.491423
$ws.Range("I9").Value = 0.04994691167843914
$ws.Range("J9").Value = 0.04994691167843914
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 139.6613298232973
$ws.Range("R9").Value = 1256.951968409676
$ws.Range("S9").Value = 0.04585653755017752
$ws.Range("T9").Value = 0.04585653755017753
$ws.Range("G10").Value = 15.497141
$ws.Range("H10").Value = 46.491423
$ws.Range("I10").Value = 0.04994691167843914
$ws.Range("J10").Value = 0.04994691167843914
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 10.744462300979
$ws.Range("R10").Value = 96.70016070881098
$ws.Range("S10").Value = 0.003527847254388099
$ws.Range("T10").Value = 0.0035278472543881
